$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (2022, month 2): mean_Y and max_Y updated
$ws.Range("C13").Value = 0.00116666666666667
$ws.Range("D13").Value = 0.121

# Row 14 (2022, month 11): mean_Y and max_Y updated
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0

# Row 15 (2022, month 12): mean_Y updated
$ws.Range("C15").Value = 0.211356182795699
